$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Rewrite ", och dessutom förklara v" -> ". Jag kommer dessutom att
#    förklara v", ending up as four separate runs (matching how Word keeps
#    runs split across distinct edit touches even when the resulting
#    character formatting is identical), without disturbing the neighboring
#    runs on either side.
# ---------------------------------------------------------------------------

# Locate the phrase to rewrite.
$target = $d.Content
$target.Find.Execute(", och dessutom förklara v", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
$base = $target.Start

# Isolate the target span from its left/right neighbors first (turn Bold on
# and keep it on through the edit) so the subsequent text rewrite doesn't
# re-coalesce it with the neighboring, identically-formatted runs.
$target.Font.Bold = 1

# Replace with a same-shaped placeholder so each of the four future pieces
# can be re-located unambiguously afterwards (avoids relying on offsets that
# could shift due to the Bold call above).
$target2 = $d.Range($base, $base + 25)
$target2.Text = "AAAAAAAAAAAABBBBBBBBBBCCCCDDDDDDDDDD"

# Write each of the four final pieces into place.
$rD = $d.Range($base + 26, $base + 36)
$rD.Text = "förklara v"

$rC = $d.Range($base + 22, $base + 26)
$rC.Text = "att "

$rB = $d.Range($base + 12, $base + 22)
$rB.Text = " dessutom "

$rA = $d.Range($base, $base + 12)
$rA.Text = ". Jag kommer"

# Clear the Bold flag piece by piece (this both restores the original look
# and keeps the four pieces as separate <w:r> elements, rather than folding
# them back into a single run).
$d.Range($base, $base + 12).Font.Bold = 0
$d.Range($base + 12, $base + 22).Font.Bold = 0
$d.Range($base + 22, $base + 26).Font.Bold = 0
$d.Range($base + 26, $base + 36).Font.Bold = 0

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark: delete it from around "sp_return_movie"
#    and drop it right before "förklara v" (between "att " and "förklara"),
#    mirroring Word's own behaviour of relocating its last-edit-position
#    bookmark to the newest edit.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$bmPoint = $d.Range($base + 26, $base + 26)
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null
